# Add the 2024/10/29 column (AY) of data to the "合成確率" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- column width for the new column AY (51) ---
# Raw OOXML <col width="12"/> corresponds to a COM ColumnWidth of 11.17
# for this sheet's default font (observed from the existing width=12 columns).
$ws.Columns.Item(51).ColumnWidth = 11.17

# --- header cell AY1: literal text "2024/10/29" (must stay text, not become a date) ---
$headerCell = $ws.Cells.Item(1, 51)
$headerCell.Value = "'2024/10/29"
$fmtSrcHeader = $ws.Cells.Item(1, 50)   # AX1, style s="1"
$fmtSrcHeader.Copy()
$headerCell.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- data rows 2-53 for column AY ---
# Reference cells carrying each of the three direct-format styles used on
# this sheet (s=1 plain, s=2 yellow fill, s=3 blue fill).
$fmtStyle1 = $ws.Cells.Item(2, 1)   # A2  -> s="1"
$fmtStyle2 = $ws.Cells.Item(2, 4)   # D2  -> s="2"
$fmtStyle3 = $ws.Cells.Item(2, 14)  # N2  -> s="3"

$values = @(120.2,198.2,197.5,111.5,117.2,164.4,181.3,238.2,151.1,162.2,169.7,168.7,191.5,187.1,144.9,109.1,156.8,158.9,155,130.3,152.3,152.9,140.4,209.5,136.8,218.3,167.8,145.6,124.2,191.4,290.5,150.6,322.8,206.9,128.5,120.7,148.6,202.8,143.1,118.2,185,149.6,114.2,172.5,155.3,143.4,145.3,116.4,203.2,171.3,190.9,117.4)
$styles = @(2,1,1,2,2,1,1,1,1,1,1,1,1,1,1,2,1,1,1,3,1,1,1,1,3,1,1,1,2,1,1,1,1,1,3,2,1,1,1,2,1,1,2,1,1,1,1,2,1,1,1,2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $rowNum = $i + 2
    $cell = $ws.Cells.Item($rowNum, 51)
    $cell.Value = $values[$i]

    $styleId = $styles[$i]
    if ($styleId -eq 2) {
        $src = $fmtStyle2
    } elseif ($styleId -eq 3) {
        $src = $fmtStyle3
    } else {
        $src = $fmtStyle1
    }
    $src.Copy()
    $cell.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
